$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 14:32:38"
$ws1.Range("A3").Value = "Total filas: 209"

# Row 56
$ws1.Range("A56").Value = "07:38:39"
$ws1.Range("B56").Value = "09:17"
$ws1.Range("C56").Value = "14_ABASTO"
$ws1.Range("D56").Value = 99
$ws1.Range("E56").Value = "LP1912"

# Row 57
$ws1.Range("A57").Value = "08:27:16"
$ws1.Range("B57").Value = "09:17"
$ws1.Range("C57").Value = "27_EL RETIRO"
$ws1.Range("D57").Value = 50
$ws1.Range("E57").Value = "LP1912"

# Row 106
$ws1.Range("A106").Value = "11:52:01"
$ws1.Range("B106").Value = "11:52"
$ws1.Range("C106").Value = "15X38_ABASTO"
$ws1.Range("D106").Value = 0
$ws1.Range("E106").Value = "LP1912"

# Row 107
$ws1.Range("A107").Value = "10:05:51"
$ws1.Range("B107").Value = "11:52"
$ws1.Range("C107").Value = "225_GOMEZ"
$ws1.Range("D107").Value = 107
$ws1.Range("E107").Value = "LP1912"

# Row 109
$ws1.Range("A109").Value = "10:37:52"
$ws1.Range("B109").Value = "11:53"
$ws1.Range("C109").Value = "23_HERNANDEZ"
$ws1.Range("D109").Value = 76
$ws1.Range("E109").Value = "LP1912"

# Row 110
$ws1.Range("A110").Value = "10:50:41"
$ws1.Range("B110").Value = "11:53"
$ws1.Range("C110").Value = "225_GOMEZ"
$ws1.Range("D110").Value = 63
$ws1.Range("E110").Value = "LP1912"

# Row 135
$ws1.Range("A135").Value = "11:34:59"
$ws1.Range("B135").Value = "12:36"
$ws1.Range("C135").Value = "23_HERNANDEZ"
$ws1.Range("D135").Value = 62
$ws1.Range("E135").Value = "LP1912"

# Row 136
$ws1.Range("A136").Value = "10:50:41"
$ws1.Range("B136").Value = "12:36"
$ws1.Range("C136").Value = "27_EL RETIRO"
$ws1.Range("D136").Value = 106
$ws1.Range("E136").Value = "LP1912"

# Row 197
$ws1.Range("A197").Value = "14:32:38"
$ws1.Range("B197").Value = "15:16"
$ws1.Range("C197").Value = "16_SANTA ANA"
$ws1.Range("D197").Value = 44
$ws1.Range("E197").Value = "LP1912"

# Row 198
$ws1.Range("A198").Value = "14:12:26"
$ws1.Range("B198").Value = "15:17"
$ws1.Range("C198").Value = "14_ABASTO"
$ws1.Range("D198").Value = 65
$ws1.Range("E198").Value = "LP1912"

# Row 199
$ws1.Range("A199").Value = "13:56:11"
$ws1.Range("B199").Value = "15:18"
$ws1.Range("C199").Value = "14_ABASTO"
$ws1.Range("D199").Value = 82
$ws1.Range("E199").Value = "LP1912"

# Row 200
$ws1.Range("A200").Value = "14:12:26"
$ws1.Range("B200").Value = "15:29"
$ws1.Range("C200").Value = "10_OLMOS"
$ws1.Range("D200").Value = 77
$ws1.Range("E200").Value = "LP1912"

# Row 201
$ws1.Range("A201").Value = "13:41:54"
$ws1.Range("B201").Value = "15:32"
$ws1.Range("C201").Value = "11_ETCHEVERRY"
$ws1.Range("D201").Value = 111
$ws1.Range("E201").Value = "LP1912"

# Row 202
$ws1.Range("A202").Value = "13:41:54"
$ws1.Range("B202").Value = "15:33"
$ws1.Range("C202").Value = "215C_EL PATO"
$ws1.Range("D202").Value = 112
$ws1.Range("E202").Value = "LP1912"

# Row 203
$ws1.Range("A203").Value = "13:56:11"
$ws1.Range("B203").Value = "15:34"
$ws1.Range("C203").Value = "215C_EL PATO"
$ws1.Range("D203").Value = 98
$ws1.Range("E203").Value = "LP1912"

# Row 204
$ws1.Range("A204").Value = "14:12:26"
$ws1.Range("B204").Value = "15:36"
$ws1.Range("C204").Value = "23_HERNANDEZ"
$ws1.Range("D204").Value = 84
$ws1.Range("E204").Value = "LP1912"

# Row 205
$ws1.Range("A205").Value = "14:12:26"
$ws1.Range("B205").Value = "15:41"
$ws1.Range("C205").Value = "11_ETCHEVERRY"
$ws1.Range("D205").Value = 89
$ws1.Range("E205").Value = "LP1912"

# Row 206
$ws1.Range("A206").Value = "13:56:11"
$ws1.Range("B206").Value = "15:42"
$ws1.Range("C206").Value = "11_ETCHEVERRY"
$ws1.Range("D206").Value = 106
$ws1.Range("E206").Value = "LP1912"

# Row 207
$ws1.Range("A207").Value = "13:56:11"
$ws1.Range("B207").Value = "15:53"
$ws1.Range("C207").Value = "15X38_ABASTO"
$ws1.Range("D207").Value = 117
$ws1.Range("E207").Value = "LP1912"

# Row 208
$ws1.Range("A208").Value = "13:56:11"
$ws1.Range("B208").Value = "15:53"
$ws1.Range("C208").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D208").Value = 117
$ws1.Range("E208").Value = "LP1912"

# Row 209
$ws1.Range("A209").Value = "14:32:38"
$ws1.Range("B209").Value = "15:55"
$ws1.Range("C209").Value = "27_EL RETIRO"
$ws1.Range("D209").Value = 83
$ws1.Range("E209").Value = "LP1912"

# Row 210
$ws1.Range("A210").Value = "14:12:26"
$ws1.Range("B210").Value = "15:56"
$ws1.Range("C210").Value = "27_EL RETIRO"
$ws1.Range("D210").Value = 104
$ws1.Range("E210").Value = "LP1912"

# Row 211
$ws1.Range("A211").Value = "14:12:26"
$ws1.Range("B211").Value = "16:05"
$ws1.Range("C211").Value = "14_ABASTO"
$ws1.Range("D211").Value = 113
$ws1.Range("E211").Value = "LP1912"

# Row 212
$ws1.Range("A212").Value = "14:32:38"
$ws1.Range("B212").Value = "16:14"
$ws1.Range("C212").Value = "17_ROMERO"
$ws1.Range("D212").Value = 102
$ws1.Range("E212").Value = "LP1912"

# Row 213
$ws1.Range("A213").Value = "14:32:38"
$ws1.Range("B213").Value = "16:17"
$ws1.Range("C213").Value = "10_OLMOS"
$ws1.Range("D213").Value = 105
$ws1.Range("E213").Value = "LP1912"

# Row 214
$ws1.Range("A214").Value = "14:32:38"
$ws1.Range("B214").Value = "16:21"
$ws1.Range("C214").Value = "23_HERNANDEZ"
$ws1.Range("D214").Value = 109
$ws1.Range("E214").Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:32:38"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 14:32:38"
$ws3.Range("A3").Value = "Total filas: 31"

# Row 34
$ws3.Range("A34").Value = "14:32:38"
$ws3.Range("B34").Value = "16:01"
$ws3.Range("C34").Value = "215C_LA PLATA"
$ws3.Range("D34").Value = 89
$ws3.Range("E34").Value = "L6203"

# Row 35
$ws3.Range("A35").Value = "14:12:26"
$ws3.Range("B35").Value = "16:02"
$ws3.Range("C35").Value = "215C_LA PLATA"
$ws3.Range("D35").Value = 110
$ws3.Range("E35").Value = "L6203"

# Row 36
$ws3.Range("A36").Value = "14:32:38"
$ws3.Range("B36").Value = "16:29"
$ws3.Range("C36").Value = "215B_LP-P MOR-40 Y 115"
$ws3.Range("D36").Value = 117
$ws3.Range("E36").Value = "L6173"
